$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("E7").Value = 16.695
$ws.Range("B9").Value = 6.439
$ws.Range("E12").Value = 17.646
$ws.Range("B13").Value = 5.302999999999999
$ws.Range("E14").Value = 16.88
$ws.Range("B16").Value = 5.294
$ws.Range("B18").Value = 5.328000000000001
$ws.Range("E19").Value = 16.566
$ws.Range("B20").Value = 6.239999999999999
$ws.Range("B26").Value = 5.897
$ws.Range("E26").Value = 16.473
$ws.Range("B27").Value = 5.752000000000001
$ws.Range("E27").Value = 16.511
$ws.Range("B29").Value = 5.532
$ws.Range("E29").Value = 17.289
$ws.Range("B35").Value = 7.657999999999999
$ws.Range("B36").Value = 7.923
$ws.Range("E37").Value = 16.855
$ws.Range("E38").Value = 16.423
$ws.Range("B45").Value = 5.516
$ws.Range("E47").Value = 16.702
$ws.Range("E51").Value = 16.769
$ws.Range("E52").Value = 16.714
$ws.Range("B55").Value = 4.711
$ws.Range("E55").Value = 16.354
$ws.Range("B57").Value = 5.147
$ws.Range("B69").Value = 5.363
$ws.Range("E69").Value = 17.321
$ws.Range("E70").Value = 17.48
$ws.Range("B76").Value = 6.425999999999999
$ws.Range("E76").Value = 16.216
$ws.Range("B78").Value = 8.254
$ws.Range("E81").Value = 16.5
$ws.Range("B82").Value = 5.286
$ws.Range("B83").Value = 5.88
$ws.Range("E83").Value = 16.917
$ws.Range("B93").Value = 5.795999999999999
$ws.Range("E94").Value = 18.151
$ws.Range("B97").Value = 5.247
$ws.Range("E100").Value = 16.527
$ws.Range("E102").Value = 16.724
